# Applies the Dec 1 2024 cryptos-list data refresh to Sheet1.
# Column D ("Price") values are prefixed with a literal apostrophe so Excel
# stores them as text (matching the original inlineStr cells) instead of
# auto-converting number-like strings (e.g. "6.76") into numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'96.390.49"
$ws.Range("E2").Value = "  -0.31%  "
$ws.Range("D3").Value = "'3.694.98"
$ws.Range("E3").Value = "  +1.30%  "
$ws.Range("D5").Value = "'236.44"
$ws.Range("E5").Value = "  -2.57%  "
$ws.Range("D6").Value = "'1.88"
$ws.Range("E6").Value = "  +2.06%  "
$ws.Range("D7").Value = "'651.12"
$ws.Range("E7").Value = "  -0.70%  "
$ws.Range("E8").Value = "  +1.11%  "
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("E10").Value = "  -1.43%  "
$ws.Range("D11").Value = "'3.695.32"
$ws.Range("E11").Value = "  +1.44%  "
$ws.Range("D12").Value = "'44.28"
$ws.Range("E12").Value = "  -0.25%  "
$ws.Range("E13").Value = "  -0.04%  "
$ws.Range("B14").Value = "Toncoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D14").Value = "'6.76"
$ws.Range("E14").Value = "  +3.77%  "
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "'0.0000293"
$ws.Range("E15").Value = "  +13.09%  "
$ws.Range("D16").Value = "'4.382.87"
$ws.Range("E16").Value = "  +1.30%  "
$ws.Range("D17").Value = "'96.213.79"
$ws.Range("E17").Value = "  -0.35%  "
$ws.Range("D18").Value = "'8.81"
$ws.Range("E18").Value = "  +13.59%  "
$ws.Range("D19").Value = "'3.698.89"
$ws.Range("E19").Value = "  +1.41%  "
$ws.Range("D20").Value = "'12.96"
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("D21").Value = "'18.83"
$ws.Range("E21").Value = "  +2.48%  "
$ws.Range("D22").Value = "'0.504"
$ws.Range("E22").Value = "  -5.79%  "
$ws.Range("D23").Value = "'517.55"
$ws.Range("E23").Value = "  +0.90%  "
$ws.Range("D24").Value = "'3.37"
$ws.Range("E24").Value = "  -2.37%  "
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("D26").Value = "'6.99"
$ws.Range("E26").Value = "  +1.11%  "
$ws.Range("D27").Value = "'100.79"
$ws.Range("E27").Value = "  -0.36%  "
$ws.Range("D28").Value = "'13.14"
$ws.Range("E28").Value = "  +0.22%  "
$ws.Range("E29").Value = "  +2.79%  "
$ws.Range("E30").Value = "  -0.94%  "
$ws.Range("D31").Value = "'12.06"
$ws.Range("E31").Value = "  +1.39%  "
$ws.Range("E32").Value = "  +0.17%  "
$ws.Range("D33").Value = "'1.85"
$ws.Range("E33").Value = "  +5.30%  "
$ws.Range("E34").Value = "  -0.83%  "
$ws.Range("D35").Value = "'0.998"
$ws.Range("E35").Value = "  -0.37%  "
$ws.Range("D36").Value = "'657.83"
$ws.Range("E36").Value = "  +6.71%  "
$ws.Range("D37").Value = "'32.14"
$ws.Range("D38").Value = "'0.586"
$ws.Range("E38").Value = "  -0.05%  "
$ws.Range("D39").Value = "'8.84"
$ws.Range("E39").Value = "  -0.08%  "
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("D41").Value = "'2.10"
$ws.Range("E41").Value = "  +8.01%  "
$ws.Range("D42").Value = "'6.84"
$ws.Range("E42").Value = "  +11.43%  "
$ws.Range("D43").Value = "'41.15"
$ws.Range("E43").Value = "  -4.24%  "
$ws.Range("E44").Value = "  +0.90%  "
$ws.Range("D45").Value = "'0.962"
$ws.Range("E45").Value = "  +0.81%  "
$ws.Range("D46").Value = "'0.0447"
$ws.Range("E46").Value = "  +0.94%  "
$ws.Range("D47").Value = "'0.436"
$ws.Range("E47").Value = "  +5.64%  "
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").Value = "'2.27"
$ws.Range("E48").Value = "  -1.89%  "
$ws.Range("B49").Value = "WhiteBITCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D49").Value = "'23.56"
$ws.Range("E49").Value = "  -0.12%  "
$ws.Range("D50").Value = "'8.46"
$ws.Range("E50").Value = "  -1.89%  "
$ws.Range("D51").Value = "'3.52"
$ws.Range("E51").Value = "  +2.41%  "
